$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "            +-------------------+                     ",
    "            |   STATE_IDLE      |<-------------------+",
    "            |  (Wait for frame) |                   |",
    "            +-------------------+                   |",
    "                     |                               |",
    "                     | MRxDV = 1 && RxStartFrm = 1   |",
    "                     v                               |",
    "            +-------------------+                   |",
    "            |    STATE_SFD      |                   |",
    "            |   (Start of frame)|                   |",
    "            +-------------------+                   |",
    "                     | MRxDV = 1                    |",
    "                     v                               |",
    "            +-------------------+                   |",
    "            |   STATE_HEADER     |                   |",
    "            | (Process Header)   |                   |",
    "            +-------------------+                   |",
    "                     | header_byte_cnt == 14         |",
    "                     v                               |",
    "            +-------------------+                   |",
    "            |    STATE_DATA      |                   |",
    "            |  (Capture Payload) |                   |",
    "            +-------------------+                   |",
    "                     |                               |",
    "                     | RxEndFrm = 1                  |",
    "                     v                               |",
    "            +-------------------+                   |",
    "            |   STATE_IDLE      |-------------------+"
)

$startRow = 25
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 4).Value = $values[$i]
}

$ws.Columns.Item(4).ColumnWidth = 39.28515625

$ws.Range("E31").Select()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
